{"js": "// The notecard's \"Obtaining the Sample Slope\" heading becomes\n// \"Obtaining the Sample (Observed) Slope\".\nconst body = context.document.body;\n\nconst results = body.search(\"Obtaining the Sample Slope\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the whole heading text in place so its formatting (bold, size 28)\n  // is preserved; the run ends up reading \"Obtaining the Sample (Observed) Slope\".\n  results.items[0].insertText(\"Obtaining the Sample (Observed) Slope\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The notecard's \"Obtaining the Sample Slope\" heading becomes\n# \"Obtaining the Sample (Observed) Slope\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Obtaining the Sample Slope\"\n$find.Replacement.Text = \"Obtaining the Sample (Observed) Slope\"\n\n# wdFindContinue = 1, wdReplaceAll = 2 -- replace across the whole document body\n# (the heading only occurs once, but this mirrors how Word's Find/Replace is\n# normally driven from PowerShell/VBA).\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
